# Adds 35 new daily rows (2021-08-26 .. 2021-09-29) to the regional COVID data
# table on Sheet1, extending it from row 358 to row 393, and updates the
# window's scroll/selection state to match, mirroring the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 359
$lastNewRow  = 393
$lastOldRow  = 358

# Propagate the date-column formatting (style/number format) from the last
# existing row ("A358") down across all the newly added rows.
$ws.Range("A$lastOldRow").Copy() | Out-Null
$ws.Range("A${firstNewRow}:A${lastNewRow}").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Write the new data (date-serial in column A, values in B:M) in one shot.
$data = New-Object 'object[,]' 35,13
$data[0,0] = 44434
$data[0,1] = 16
$data[0,2] = 14
$data[0,3] = 5
$data[0,4] = 3
$data[0,5] = 11
$data[0,6] = 16
$data[0,7] = 1
$data[0,8] = 10
$data[0,9] = 5
$data[0,10] = 8
$data[0,11] = 5
$data[0,12] = 3
$data[1,0] = 44435
$data[1,1] = 12
$data[1,2] = 14
$data[1,3] = 2
$data[1,4] = 6
$data[1,5] = 14
$data[1,6] = 9
$data[1,7] = 1
$data[1,8] = 10
$data[1,9] = 10
$data[1,10] = 7
$data[1,11] = 3
$data[1,12] = 4
$data[2,0] = 44436
$data[2,1] = 17
$data[2,2] = 7
$data[2,3] = 6
$data[2,4] = 7
$data[2,5] = 12
$data[2,6] = 11
$data[2,7] = 0
$data[2,8] = 9
$data[2,9] = 9
$data[2,10] = 10
$data[2,11] = 1
$data[2,12] = 4
$data[3,0] = 44437
$data[3,1] = 12
$data[3,2] = 2
$data[3,3] = 2
$data[3,4] = 5
$data[3,5] = 11
$data[3,6] = 11
$data[3,7] = 0
$data[3,8] = 6
$data[3,9] = 11
$data[3,10] = 13
$data[3,11] = 1
$data[3,12] = 2
$data[4,0] = 44438
$data[4,1] = 20
$data[4,2] = 17
$data[4,3] = 7
$data[4,4] = 7
$data[4,5] = 14
$data[4,6] = 9
$data[4,7] = 0
$data[4,8] = 9
$data[4,9] = 7
$data[4,10] = 8
$data[4,11] = 2
$data[4,12] = 3
$data[5,0] = 44439
$data[5,1] = 25
$data[5,2] = 15
$data[5,3] = 4
$data[5,4] = 7
$data[5,5] = 21
$data[5,6] = 12
$data[5,7] = 0
$data[5,8] = 6
$data[5,9] = 7
$data[5,10] = 6
$data[5,11] = 4
$data[5,12] = 2
$data[6,0] = 44440
$data[6,1] = 21
$data[6,2] = 7
$data[6,3] = 5
$data[6,4] = 7
$data[6,5] = 19
$data[6,6] = 12
$data[6,7] = 0
$data[6,8] = 9
$data[6,9] = 11
$data[6,10] = 7
$data[6,11] = 0
$data[6,12] = 2
$data[7,0] = 44441
$data[7,1] = 17
$data[7,2] = 5
$data[7,3] = 3
$data[7,4] = 7
$data[7,5] = 4
$data[7,6] = 13
$data[7,7] = 0
$data[7,8] = 6
$data[7,9] = 5
$data[7,10] = 7
$data[7,11] = 1
$data[7,12] = 2
$data[8,0] = 44442
$data[8,1] = 23
$data[8,2] = 12
$data[8,3] = 6
$data[8,4] = 5
$data[8,5] = 12
$data[8,6] = 15
$data[8,7] = 0
$data[8,8] = 4
$data[8,9] = 13
$data[8,10] = 9
$data[8,11] = 1
$data[8,12] = 4
$data[9,0] = 44443
$data[9,1] = 16
$data[9,2] = 6
$data[9,3] = 1
$data[9,4] = 4
$data[9,5] = 7
$data[9,6] = 13
$data[9,7] = 0
$data[9,8] = 2
$data[9,9] = 9
$data[9,10] = 7
$data[9,11] = 2
$data[9,12] = 3
$data[10,0] = 44444
$data[10,1] = 18
$data[10,2] = 11
$data[10,3] = 4
$data[10,4] = 2
$data[10,5] = 10
$data[10,6] = 11
$data[10,7] = 1
$data[10,8] = 7
$data[10,9] = 8
$data[10,10] = 6
$data[10,11] = 2
$data[10,12] = 0
$data[11,0] = 44445
$data[11,1] = 15
$data[11,2] = 5
$data[11,3] = 3
$data[11,4] = 4
$data[11,5] = 12
$data[11,6] = 4
$data[11,7] = 1
$data[11,8] = 6
$data[11,9] = 9
$data[11,10] = 8
$data[11,11] = 4
$data[11,12] = 1
$data[12,0] = 44446
$data[12,1] = 16
$data[12,2] = 7
$data[12,3] = 2
$data[12,4] = 4
$data[12,5] = 17
$data[12,6] = 6
$data[12,7] = 0
$data[12,8] = 8
$data[12,9] = 8
$data[12,10] = 9
$data[12,11] = 0
$data[12,12] = 2
$data[13,0] = 44447
$data[13,1] = 14
$data[13,2] = 15
$data[13,3] = 2
$data[13,4] = 2
$data[13,5] = 6
$data[13,6] = 4
$data[13,7] = 0
$data[13,8] = 8
$data[13,9] = 8
$data[13,10] = 7
$data[13,11] = 4
$data[13,12] = 2
$data[14,0] = 44448
$data[14,1] = 20
$data[14,2] = 5
$data[14,3] = 2
$data[14,4] = 3
$data[14,5] = 14
$data[14,6] = 6
$data[14,7] = 1
$data[14,8] = 3
$data[14,9] = 9
$data[14,10] = 7
$data[14,11] = 1
$data[14,12] = 3
$data[15,0] = 44449
$data[15,1] = 14
$data[15,2] = 2
$data[15,3] = 4
$data[15,4] = 3
$data[15,5] = 16
$data[15,6] = 6
$data[15,7] = 0
$data[15,8] = 6
$data[15,9] = 9
$data[15,10] = 2
$data[15,11] = 3
$data[15,12] = 1
$data[16,0] = 44450
$data[16,1] = 14
$data[16,2] = 5
$data[16,3] = 1
$data[16,4] = 4
$data[16,5] = 4
$data[16,6] = 9
$data[16,7] = 0
$data[16,8] = 4
$data[16,9] = 4
$data[16,10] = 3
$data[16,11] = 2
$data[16,12] = 2
$data[17,0] = 44451
$data[17,1] = 16
$data[17,2] = 4
$data[17,3] = 2
$data[17,4] = 5
$data[17,5] = 13
$data[17,6] = 3
$data[17,7] = 1
$data[17,8] = 1
$data[17,9] = 6
$data[17,10] = 6
$data[17,11] = 0
$data[17,12] = 1
$data[18,0] = 44452
$data[18,1] = 15
$data[18,2] = 6
$data[18,3] = 3
$data[18,4] = 4
$data[18,5] = 15
$data[18,6] = 9
$data[18,7] = 1
$data[18,8] = 2
$data[18,9] = 7
$data[18,10] = 9
$data[18,11] = 0
$data[18,12] = 1
$data[19,0] = 44453
$data[19,1] = 16
$data[19,2] = 5
$data[19,3] = 1
$data[19,4] = 4
$data[19,5] = 12
$data[19,6] = 4
$data[19,7] = 0
$data[19,8] = 2
$data[19,9] = 6
$data[19,10] = 12
$data[19,11] = 0
$data[19,12] = 3
$data[20,0] = 44454
$data[20,1] = 7
$data[20,2] = 5
$data[20,3] = 2
$data[20,4] = 3
$data[20,5] = 14
$data[20,6] = 0
$data[20,7] = 0
$data[20,8] = 2
$data[20,9] = 3
$data[20,10] = 8
$data[20,11] = 1
$data[20,12] = 1
$data[21,0] = 44455
$data[21,1] = 10
$data[21,2] = 3
$data[21,3] = 2
$data[21,4] = 2
$data[21,5] = 11
$data[21,6] = 5
$data[21,7] = 0
$data[21,8] = 5
$data[21,9] = 3
$data[21,10] = 4
$data[21,11] = 0
$data[21,12] = 1
$data[22,0] = 44456
$data[22,1] = 7
$data[22,2] = 5
$data[22,3] = 1
$data[22,4] = 3
$data[22,5] = 10
$data[22,6] = 9
$data[22,7] = 0
$data[22,8] = 3
$data[22,9] = 5
$data[22,10] = 6
$data[22,11] = 1
$data[22,12] = 3
$data[23,0] = 44457
$data[23,1] = 14
$data[23,2] = 4
$data[23,3] = 3
$data[23,4] = 5
$data[23,5] = 10
$data[23,6] = 0
$data[23,7] = 0
$data[23,8] = 1
$data[23,9] = 4
$data[23,10] = 7
$data[23,11] = 0
$data[23,12] = 0
$data[24,0] = 44458
$data[24,1] = 8
$data[24,2] = 3
$data[24,3] = 1
$data[24,4] = 3
$data[24,5] = 8
$data[24,6] = 2
$data[24,7] = 0
$data[24,8] = 2
$data[24,9] = 2
$data[24,10] = 4
$data[24,11] = 0
$data[24,12] = 1
$data[25,0] = 44459
$data[25,1] = 11
$data[25,2] = 2
$data[25,3] = 1
$data[25,4] = 2
$data[25,5] = 19
$data[25,6] = 5
$data[25,7] = 0
$data[25,8] = 2
$data[25,9] = 1
$data[25,10] = 5
$data[25,11] = 0
$data[25,12] = 0
$data[26,0] = 44460
$data[26,1] = 10
$data[26,2] = 2
$data[26,3] = 3
$data[26,4] = 5
$data[26,5] = 16
$data[26,6] = 4
$data[26,7] = 0
$data[26,8] = 3
$data[26,9] = 5
$data[26,10] = 1
$data[26,11] = 0
$data[26,12] = 2
$data[27,0] = 44461
$data[27,1] = 8
$data[27,2] = 4
$data[27,3] = 0
$data[27,4] = 2
$data[27,5] = 9
$data[27,6] = 1
$data[27,7] = 0
$data[27,8] = 1
$data[27,9] = 0
$data[27,10] = 6
$data[27,11] = 0
$data[27,12] = 0
$data[28,0] = 44462
$data[28,1] = 11
$data[28,2] = 3
$data[28,3] = 2
$data[28,4] = 2
$data[28,5] = 9
$data[28,6] = 1
$data[28,7] = 0
$data[28,8] = 3
$data[28,9] = 1
$data[28,10] = 3
$data[28,11] = 0
$data[28,12] = 1
$data[29,0] = 44463
$data[29,1] = 6
$data[29,2] = 3
$data[29,3] = 1
$data[29,4] = 2
$data[29,5] = 6
$data[29,6] = 4
$data[29,7] = 0
$data[29,8] = 3
$data[29,9] = 1
$data[29,10] = 0
$data[29,11] = 0
$data[29,12] = 0
$data[30,0] = 44464
$data[30,1] = 8
$data[30,2] = 2
$data[30,3] = 0
$data[30,4] = 2
$data[30,5] = 5
$data[30,6] = 3
$data[30,7] = 0
$data[30,8] = 3
$data[30,9] = 5
$data[30,10] = 1
$data[30,11] = 0
$data[30,12] = 1
$data[31,0] = 44465
$data[31,1] = 5
$data[31,2] = 6
$data[31,3] = 0
$data[31,4] = 1
$data[31,5] = 10
$data[31,6] = 3
$data[31,7] = 0
$data[31,8] = 3
$data[31,9] = 2
$data[31,10] = 4
$data[31,11] = 0
$data[31,12] = 1
$data[32,0] = 44466
$data[32,1] = 8
$data[32,2] = 3
$data[32,3] = 0
$data[32,4] = 2
$data[32,5] = 16
$data[32,6] = 2
$data[32,7] = 0
$data[32,8] = 0
$data[32,9] = 0
$data[32,10] = 1
$data[32,11] = 0
$data[32,12] = 0
$data[33,0] = 44467
$data[33,1] = 5
$data[33,2] = 3
$data[33,3] = 0
$data[33,4] = 2
$data[33,5] = 6
$data[33,6] = 3
$data[33,7] = 0
$data[33,8] = 1
$data[33,9] = 1
$data[33,10] = 2
$data[33,11] = 2
$data[33,12] = 1
$data[34,0] = 44468
$data[34,1] = 5
$data[34,2] = 3
$data[34,3] = 1
$data[34,4] = 3
$data[34,5] = 8
$data[34,6] = 2
$data[34,7] = 0
$data[34,8] = 1
$data[34,9] = 0
$data[34,10] = 0
$data[34,11] = 0
$data[34,12] = 0

$ws.Range("A${firstNewRow}:M${lastNewRow}").Value = $data

# Match the saved view state: scrolled down near the bottom of the sheet
# with a single cell selected just below the new data.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 352
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B357").Select()
